$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.204.41"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "3.825.18"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "446.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.90%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.744"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000319"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.70%  "

$ws.Range("D14").Value = "4.398.22"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.47%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.868.00"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.46%  "

$ws.Range("E19").Value = "  +7.76%  "

$ws.Range("D20").Value = "67.139.93"
$ws.Range("E20").Value = "  +0.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "421.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.77%  "

$ws.Range("E23").Value = "  +9.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.52%  "

$ws.Range("B25").Value = "EthereumClassic"
$ws.Range("C25").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "37.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.28%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +26.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "732.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.66%  "

$ws.Range("E31").Value = "  +12.32%  "

$ws.Range("E32").Value = "  +10.79%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "44.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +18.40%  "

$ws.Range("E35").Value = "  +7.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "56.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.84%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +24.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0480"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.81%  "

$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0684"
$ws.Range("E41").Value = "  -10.64%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.141"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.55%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.336"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +22.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "144.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.02%  "

